$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" (Changed) date column C for all data rows (2-89)
# from 45233 (2023-11-03) to 45243 (2023-11-13).
$ws.Range("C2:C89").Value = 45243
